$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = 9.000
$ws.Range("H12").Value = 851.89

$ws.Range("F19").Value = 5.000
$ws.Range("H19").Value = 279.53

$ws.Range("H24").Value = 9497.27

$ws.Range("F28").Value = 1196.000
$ws.Range("H28").Value = 10853.42

$ws.Range("F30").Value = 1358.000
$ws.Range("H30").Value = 13955.51

$ws.Range("H45").Value = 37404.19

$ws.Range("H49").Value = 25874.67

$ws.Range("F57").Value = 4085.000
$ws.Range("H57").Value = 15872.05

$ws.Range("F76").Value = 5.000
$ws.Range("H76").Value = 1368.67

$ws.Range("H88").Value = 14137.37

$ws.Range("F89").Value = 23654.000
$ws.Range("H89").Value = 35963.12

$ws.Range("F108").Value = 5.000
$ws.Range("H108").Value = 695.82

$ws.Range("F117").Value = 153.820
$ws.Range("H117").Value = 3328.35

$ws.Range("F119").Value = 10066.300
$ws.Range("H119").Value = 22840.97

$ws.Range("F129").Value = 81.000
$ws.Range("H129").Value = 474.25
